# Update "principais_noticias" worksheet with the newest news items.
# Rows 3-11 (columns B=Title, C=Link, D=Abstract) are refreshed: the three
# oldest stories are dropped, everything shifts up, and three brand new
# stories are inserted in rows 5-7.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @{
    3 = @(
        "Allianz fecha parceria com o Rappi para distribuição de seguro via corretor",
        "https://www.revistaapolice.com.br/2023/03/allianz-fecha-parceria-com-o-rappi-para-distribuicao-de-seguro-via-corretor/",
        "Ação é voltada é exclusivamente aos clientes Rappi Prime e abrange todas as ofertas e pacotes do seguro Residencial da companhia"
    )
    4 = @(
        "Mapfre anuncia Afonso Arinos como head Comercial de Agronegócio",
        "https://www.revistaapolice.com.br/2023/03/mapfre-apresenta-afonso-arinos-como-head-comercial-de-agronegocio/",
        "Com 20 anos de experiência em posições de liderança comercial, novo executivo irá liderar estratégias comerciais para a companhia"
    )
    5 = @(
        "Brasilprev lança iniciativa de capacitação exclusiva para PCD's",
        "https://www.revistaapolice.com.br/2023/03/brasilprev-lanca-iniciativa-de-capacitacao-exclusiva-para-pcds/",
        "Com o objetivo de atrair, capacitar e contratar pessoas com deficiência, programa Lado a Lado irá treinar 50 profissionais em Tecnologia"
    )
    6 = @(
        "Metaverso em seguros: Uma realidade distante?",
        "https://www.revistaapolice.com.br/2023/03/metaverso-em-seguros-uma-realidade-distante/",
        "O Metaverso oferece um potencial significativo ao setor, atuando, por exemplo, como uma nova abordagem para alcançar o público mais jovem"
    )
    7 = @(
        "Aumento de desastres naturais reforçam importância do seguro residencial",
        "https://www.revistaapolice.com.br/2023/03/aumento-de-desastres-naturais-reforcam-importancia-do-seguro-residencial/",
        "Segundo dados da FenaSeg, apenas 15% das residências brasileiras que contam com um seguro residencial possuem cobertura contra enchentes"
    )
    8 = @(
        "Akad disponibiliza seguro por assinatura para pequenas e médias empresas",
        "https://www.revistaapolice.com.br/2023/03/akad-disponibiliza-seguro-por-assinatura-para-pequenas-e-medias-empresas/",
        "Parceria com a MarketUP e a PDV Box é aposta da seguradora para chegar a mais de 150 mil usuários ativos"
    )
    9 = @(
        "Tokio Marine registra crescimento de 40,6% na produção em 2022",
        "https://www.revistaapolice.com.br/2023/03/tokio-marine-registra-crescimento-de-406-na-producao-em-2022/",
        "A produção da companhia atingiu a marca de R$ 10,6 bilhões em faturamento, meta inicialmente prevista para ser alcançada apenas em 2024"
    )
    10 = @(
        "PDMS quer levar seguro para fatia maior da população",
        "https://www.revistaapolice.com.br/2023/03/pdms-quer-levar-seguro-para-fatia-maior-da-populacao/",
        "Plano de Desenvolvimento do Mercado de Seguros possui 65 diretrizes em 4 áreas para aumentar o nível de proteção do País"
    )
    11 = @(
        "Mapfre patrocina musical sobre a história do grupo Los Hermanos",
        "https://www.revistaapolice.com.br/2023/03/mapfre-patrocina-musical-sobre-a-historia-do-grupo-los-hermanos/",
        "Espetáculo, que tem a direção de Michel Melamed, estreia nesta sexta-feira, dia 17 de março, no Teatro Casa Grande, no Rio de Janeiro"
    )
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    $ws.Cells.Item($r, 2).Value = $vals[0]
    $ws.Cells.Item($r, 3).Value = $vals[1]
    $ws.Cells.Item($r, 4).Value = $vals[2]
}
